$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E (price & volume) columns to remain text so numeric-looking
# strings (e.g. "6.26", "142.70") are not auto-converted to numbers by Excel.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '63.868.97'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '3.054.28'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '559.45'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").Value = '142.70'
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = '3.054.38'
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("E9").Value = '  +3.19%  '
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("D11").Value = '6.26'
$ws.Range("E11").Value = '  -10.95%  '
$ws.Range("D12").Value = '0.493'
$ws.Range("E12").Value = '  +7.57%  '
$ws.Range("D13").Value = '0.0000231'
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").Value = '35.80'
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").Value = '3.554.29'
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("D16").Value = '63.950.02'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").Value = '3.058.39'
$ws.Range("E17").Value = '  -1.42%  '
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").Value = '6.81'
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("D20").Value = '476.58'
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("D21").Value = '14.10'
$ws.Range("E21").Value = '  +3.17%  '
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D22").Value = '14.65'
$ws.Range("E22").Value = '  +11.26%  '
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").Value = '0.684'
$ws.Range("E23").Value = '  +2.31%  '
$ws.Range("E24").Value = '  +0.60%  '
$ws.Range("D25").Value = '82.76'
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").Value = '2.81'
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").Value = '  +3.15%  '
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("D31").Value = '26.30'
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("D35").Value = '6.23'
$ws.Range("E35").Value = '  +2.21%  '
$ws.Range("D36").Value = '54.55'
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("D38").Value = '447.54'
$ws.Range("E38").Value = '  -2.90%  '
$ws.Range("D39").Value = '0.0815'
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").Value = '2.82'
$ws.Range("E40").Value = '  +5.37%  '
$ws.Range("D41").Value = '3.016.26'
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("D43").Value = '8.30'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").Value = '0.271'
$ws.Range("E44").Value = '  +4.92%  '
$ws.Range("D45").Value = '28.32'
$ws.Range("E45").Value = '  +1.16%  '
$ws.Range("D46").Value = '2.29'
$ws.Range("E46").Value = '  +10.49%  '
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("D49").Value = '117.70'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("E51").Value = '  +1.81%  '

# Restore the original (default) cell style now that values are set,
# so no stray number-format styling is left on the cells.
$priceVolRange.Style = "Normal"

